$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.022.94"
$ws.Range("E2").Value = "  +1.56%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.519.44"
$ws.Range("E3").Value = "  +3.02%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "535.51"
$ws.Range("E5").Value = "  +1.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.33"
$ws.Range("E6").Value = "  +2.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +2.52%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.517.91"
$ws.Range("E9").Value = "  +2.78%  "
$ws.Range("E10").Value = "  +2.80%  "
$ws.Range("E11").Value = "  -1.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.36"
$ws.Range("E12").Value = "  +1.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.348"
$ws.Range("E13").Value = "  +2.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.963.23"
$ws.Range("E14").Value = "  +2.72%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.02"
$ws.Range("E15").Value = "  +3.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "58.926.59"
$ws.Range("E16").Value = "  +1.56%  "
$ws.Range("E17").Value = "  +1.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.505.71"
$ws.Range("E18").Value = "  +2.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.06"
$ws.Range("E19").Value = "  +4.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.26"
$ws.Range("E20").Value = "  +2.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "321.72"
$ws.Range("E21").Value = "  +1.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.96"
$ws.Range("E23").Value = "  +5.71%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.11"
$ws.Range("E24").Value = "  +5.26%  "
$ws.Range("E25").Value = "  +4.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.165"
$ws.Range("E26").Value = "  +1.39%  "
$ws.Range("E27").Value = "  +1.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.51"
$ws.Range("E28").Value = "  +0.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.65"
$ws.Range("E29").Value = "  +2.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0767"
$ws.Range("E30").Value = "  +2.55%  "
$ws.Range("E31").Value = "  +0.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "170.61"
$ws.Range("E33").Value = "  +11.15%  "
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("E35").Value = "  +2.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.36"
$ws.Range("E36").Value = "  +1.54%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.06"
$ws.Range("E37").Value = "  +1.90%  "
$ws.Range("E38").Value = "  +0.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.82"
$ws.Range("E39").Value = "  +1.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.810"
$ws.Range("E40").Value = "  +4.64%  "
$ws.Range("E41").Value = "  +1.88%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "283.92"
$ws.Range("E42").Value = "  +5.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.13"
$ws.Range("E43").Value = "  +3.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  -0.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.607"
$ws.Range("E45").Value = "  +4.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "129.90"
$ws.Range("E46").Value = "  +9.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.89"
$ws.Range("E47").Value = "  +0.61%  "
$ws.Range("E48").Value = "  +0.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0503"
$ws.Range("E49").Value = "  +0.63%  "
$ws.Range("E50").Value = "  +0.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.40"
$ws.Range("E51").Value = "  +3.31%  "
